$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix a renamed event title (shared string text fix)
$ws.Range("B289").Value = "UNREAL VENDEX ALL NIGHT LONG"

# 2. Populate the new event rows 296-308.
#    Each row: seed formatting from the last populated row (295) via
#    PasteSpecial (formats only) so the new cells pick up the same cell
#    styles (s=4 for dates, s=3 for text/links) as the rest of the table,
#    then fill in values + a real hyperlink object for column E.

# 2a. Grab the blank trailing-row formatting (s=4/s=5) now, while rows
#     296-304 are still empty placeholders, and stamp it onto the newly
#     appended rows 309-331 *before* 296-304 get overwritten with data --
#     otherwise there is no more "blank" row left downstream to copy from.
$ws.Range("A296:E296").Copy()
$ws.Range("A309:E331").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A295:E295").Copy()
$ws.Range("A296:E308").PasteSpecial(-4122)
$excel.CutCopyMode = 0

function Add-EventRow {
    # NB: positional params only -- named args (-Row 296) silently fail to
    # bind in this PowerShell host, so everything below is called
    # positionally.
    param($Row, $Date, $Event, $Location, $City, $Link)

    $ws.Range("A$Row").Value = $Date
    $ws.Range("B$Row").Value = $Event
    $ws.Range("C$Row").Value = $Location
    $ws.Range("D$Row").Value = $City

    $ws.Hyperlinks.Add($ws.Range("E$Row"), $Link, "", "", $Link)

    # Re-apply the plain data-cell format to E (Hyperlinks.Add otherwise
    # bumps it to Excel's builtin blue/underline "Hyperlink" style, which
    # this sheet does not use -- its link styling lives in the shared
    # string's own run formatting instead).
    $ws.Range("E295").Copy()
    $ws.Range("E$Row").PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}

Add-EventRow 296 45759 "DISCLOSED" "PM93" "Essen" "https://www.instagram.com/reel/DIJ6ayBsB0R/?igsh=MTBtMGhlMzVtYmV0dg=="
Add-EventRow 297 45786 "X-BASS" "viersieben" "Duisburg" "https://www.instagram.com/reel/DIL43dWgM_T/?igsh=MXB5aGsyZ2JjcGh4cA=="
Add-EventRow 298 45786 "ESCALATE" "THEGATE" "Düsseldorf" "https://www.instagram.com/reel/DIMQQwooLhx/?igsh=cjVmNTZoeHltenZj"
Add-EventRow 299 45759 "SECRETRAVES" "check event link" "Köln" "https://chat.whatsapp.com/KYWhxS6wl5T6Rdgr0A5pmb"
Add-EventRow 300 45764 "SECRETRAVES" "check event link" "Köln" "https://chat.whatsapp.com/KYWhxS6wl5T6Rdgr0A5pmb"
Add-EventRow 301 45787 "AREA 15 BOCHUM & FRIENDS" "AREA 15" "Bochum" "https://www.instagram.com/s/aGlnaGxpZ2h0OjE3OTcyMzI3MTcxNzIyOTcy?story_media_id=3605001403622987852&igsh=ZndpdmFwMW5obmpo"
Add-EventRow 302 45756 "180 MIN RAVE" "Prismatic" "Dortmund" "https://www.instagram.com/reel/DIJ7zsWs8OE/?igsh=MWl1YW1tNjdteGNtNQ=="
Add-EventRow 303 45771 "LUCID VOID" "Klub Kulb" "Düsseldorf" "https://www.instagram.com/area51.techno?igsh=MWI0amhkbHZsN2RneQ=="
Add-EventRow 304 45758 "BERRYLECTRO SECRET RAVE" "check event link" "Witten" "https://chat.whatsapp.com/Ie6PzGtqJaI4oSlkXJ2Cwf?fbclid=PAZXh0bgNhZW0CMTEAAacXHqe4Af4vx6ypy4LszKZGeOxWUN3JfHG-m7ZDWcbrneqdygS1qIakMmlX-Q_aem_yJBk9a1vpWKZ4PCUJu9Oug&e=AT12jcyFqY0wes8CybzLulBZKydbvyDLlhO3n-xBWBdD3GgRQCV1EXkgR_KN_7ncwJJmyWB2VbhNNJT2ncjkmDv_VuemQwfdrQTaXHdhKg4zOmg9ZtFBcE0"
Add-EventRow 305 45758 "RAVE INDUSTRY" "Stollen134" "Dortmund" "https://www.instagram.com/reel/DG-ygzFKEW-/?igsh=MTd4b2k4aTI5ZTZrMQ=="
Add-EventRow 306 45779 "SYNDEN" "Elektroküche" "Köln" "https://www.instagram.com/_synden?igsh=MWhsZnNwdHV0d3dobw=="
Add-EventRow 307 45772 "TECHNOALLIANZ SCHRANZ ONLY" "Elektroküche" "Köln" "https://www.instagram.com/p/DH8DqrgAjRw/?igsh=ZTdpajRzaTBnZ3lu"
Add-EventRow 308 45780 "EHRENKLUB" "Odonien" "Köln" "https://t.rausgegangen.de/tickets/ehrenklub-in-odonien-20-mit"
